$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9; this shifts the existing rows 9-15 down to 10-16
# and carries the D-column date style (s="2") down with them.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with a new data record
# (same reference data as the surrounding rows, but its own measurements).
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = "1/18/2023"
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 2500
$ws.Range("O9").Value = 2500
$ws.Range("P9").Value = 2500
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Provincia de Diguillín"
$ws.Range("S9").Value = 1250
$ws.Range("T9").Value = 2

# Row 12 (old row 11, shifted down) also got its max/avg price and $/Kg
# price corrected on top of the shift.
$ws.Range("O12").Value = 2500
$ws.Range("P12").Value = 2500
$ws.Range("S12").Value = 1250
